# =====================================================================
# sync: Export all json sorted by `RowOrder`
#
# Applies:
#  1. TechTree sheet: rename "Order" column -> "RowOrder" (table + header
#     cell), widen the RowOrder column, add a green-highlight conditional
#     format over the table body.
#  2. New "TechTree_Expanded" sheet: title cell explaining the export,
#     large (22pt) font, same conditional-format highlight rule.
#  3. "ResearchRecipes" sheet: add the same conditional-format highlight
#     rule (data itself is unchanged).
#  4. New "Exploration_Order" sheet: ORDER / BiomeBlue2 / RowOrder table
#     listing the biome-unlock exploration order, backed by a
#     ListObject ("Exploration_Order_Data").
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. TechTree: rename the "Order" column to "RowOrder"
# ---------------------------------------------------------------------
$techTree = $wb.Worksheets.Item("TechTree")
$techTree.Range("R1").Value = "RowOrder"

# Widen column R (RowOrder) to fit the new header text.
$techTree.Columns.Item(18).ColumnWidth = 11.82

# Highlight rule: cells equal to "x" get dark-green text on a light-green
# fill (matches the rule added to every data sheet in this release).
$techTreeCF = $techTree.Range("A1:R240").FormatConditions.Add(1, 3, '"x"')
$techTreeCF.Font.Color = 25600
$techTreeCF.Interior.Color = 9498256

# ---------------------------------------------------------------------
# 2. New sheet: TechTree_Expanded (placed right after TechTree)
# ---------------------------------------------------------------------
$expanded = $wb.Worksheets.Add($null, $techTree)
$expanded.Name = "TechTree_Expanded"
$expanded.Range("A1").Value = "TechTree with columns expanded as multiple rows"
$expanded.Range("A1").Font.Size = 22

$expandedCF = $expanded.Range("A1").FormatConditions.Add(1, 3, '"x"')
$expandedCF.Font.Color = 25600
$expandedCF.Interior.Color = 9498256

# ---------------------------------------------------------------------
# 3. ResearchRecipes: same highlight rule, data untouched
# ---------------------------------------------------------------------
$recipes = $wb.Worksheets.Item("ResearchRecipes")
$recipesCF = $recipes.Range("A1:G10").FormatConditions.Add(1, 3, '"x"')
$recipesCF.Font.Color = 25600
$recipesCF.Interior.Color = 9498256

# ---------------------------------------------------------------------
# 4. New sheet: Exploration_Order (placed right after ResearchRecipes)
# ---------------------------------------------------------------------
$exploration = $wb.Worksheets.Add($null, $recipes)
$exploration.Name = "Exploration_Order"

$exploration.Range("A1").Value = "ORDER"
$exploration.Range("B1").Value = "BiomeBlue2"
$exploration.Range("C1").Value = "RowOrder"

$explorationData = @(
    @(0, 'BiomeBlue2_start', 0),
    @(1, 'BiomeScrapara', 1),
    @(2, 'BiomeGreen', 2),
    @(3, 'BiomeBlue2, BiomeScrapara, BiomeGreen', 3),
    @(4, 'BiomeToxicwaste', 4),
    @(5, 'BiomeConcrete', 5),
    @(6, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 6),
    @(7, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 7),
    @(8, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 8),
    @(9, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 9),
    @(10, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 10),
    @(11, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 11),
    @(12, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 12),
    @(13, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 13),
    @(14, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 14),
    @(15, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 15),
    @(16, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 16),
    @(17, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 17),
    @(18, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 18),
    @(19, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 19),
    @(20, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 20),
    @(21, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 21),
    @(22, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 22),
    @(23, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 23),
    @(24, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 24),
    @(25, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 25),
    @(26, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 26),
    @(27, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 27),
    @(28, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 28),
    @(29, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 29),
    @(30, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 30),
    @(31, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 31),
    @(32, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 32),
    @(33, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 33),
    @(34, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 34),
    @(35, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 35),
    @(36, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 36),
    @(37, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 37),
    @(38, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 38),
    @(39, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 39),
    @(40, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 40),
    @(41, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 41),
    @(42, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 42),
    @(43, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 43),
    @(44, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 44),
    @(45, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 45),
    @(46, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 46),
    @(47, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 47),
    @(48, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 48),
    @(49, 'BiomeBlue2, BiomeScrapara, BiomeGreen, BiomeToxicwaste, BiomeConcrete', 49)
)

$r = 2
foreach ($row in $explorationData) {
    $exploration.Cells.Item($r, 1).Value = $row[0]
    $exploration.Cells.Item($r, 2).Value = $row[1]
    $exploration.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

$exploration.Columns.Item(1).ColumnWidth = 8.98
$exploration.Columns.Item(2).ColumnWidth = 68.42
$exploration.Columns.Item(3).ColumnWidth = 11.82

$explorationTable = $exploration.ListObjects.Add(1, $exploration.Range("A1:C51"), $null, 1)
$explorationTable.Name = "Exploration_Order_Data"
$explorationTable.TableStyle = "TableStyleLight5"
